$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14; this pushes existing rows 14-31 down to 15-32
$ws.Rows("14:14").Insert()

# Copy the date number format from the row above (row 15, originally row 14) into new D14
$ws.Range("D15").Copy()
$ws.Range("D14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row 14 with the new data record
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "Vega Monumental Concepción"
$ws.Range("C14").Value = "Bíobío"
$ws.Range("D14").Value = 44797
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 100112026
$ws.Range("G14").Value = "Haba"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 11200
$ws.Range("N14").Value = "$/saco 25 kilos"
$ws.Range("O14").Value = "Región de Coquimbo"
$ws.Range("P14").Value = 448
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"

$wb.Save()
